$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.210.45"
$ws.Range("E2").Value = "  +6.58%  "
$ws.Range("D3").Value = "3.616.45"
$ws.Range("E3").Value = "  +9.44%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "642.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.51"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.407"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.34%  "
$ws.Range("D11").Value = "3.614.15"
$ws.Range("E11").Value = "  +9.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.202"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.92%  "
$ws.Range("D15").Value = "4.290.55"
$ws.Range("E15").Value = "  +9.29%  "
$ws.Range("D16").Value = "97.117.01"
$ws.Range("E16").Value = "  +6.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000256"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.88%  "
$ws.Range("D18").Value = "3.622.28"
$ws.Range("E18").Value = "  +9.61%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +21.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.505"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +14.51%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "517.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000201"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "98.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +20.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.146"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.184"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "30.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.575"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "584.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.95%  "
$ws.Range("E39").Value = "  +10.51%  "
$ws.Range("E40").Value = "  +4.70%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.927"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.44%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0438"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.95%  "
$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.41%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.84%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.95%  "
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.75%  "
